$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.040.68'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '2.411.55'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.71'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.67'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.506'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  +6.39%  '
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.330'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.63'
$ws.Range('E12').Value = '  -5.09%  '
$ws.Range('D13').Value = '67.932.05'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = '2.854.39'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000173'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.81'
$ws.Range('E16').Value = '  -4.13%  '
$ws.Range('D17').Value = '2.410.73'
$ws.Range('E17').Value = '  -3.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.37'
$ws.Range('E18').Value = '  -4.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '330.40'
$ws.Range('E19').Value = '  -2.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.85'
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.79'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.87'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.13'
$ws.Range('E24').Value = '  -2.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.65'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').Value = '2.533.83'
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.14'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').Value = '0.0₃0807'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.07'
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '420.00'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.60'
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '159.99'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.00'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.78'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.104'
$ws.Range('E38').Value = '  -5.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.295'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.29'
$ws.Range('E40').Value = '  -4.10%  '
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '131.84'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.31'
$ws.Range('E44').Value = '  -1.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.97'
$ws.Range('E45').Value = '  -5.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0710'
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.478'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.553'
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0913'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.37'
$ws.Range('E51').Value = '  -3.45%  '
